# Update "想去人数" (interested-count) values in column F on both the
# "展览" sheet and the "全部类型" sheet, which duplicate the same events.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value for sheet "展览" ("F" column)
$exhibitUpdates = @{
    8  = 7499
    10 = 7690
    11 = 24
    13 = 18
    14 = 6273
    15 = 3288
    27 = 3683
    32 = 1326
    35 = 2643
    36 = 1574
    37 = 16
    40 = 3358
    41 = 200
    46 = 1313
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value for sheet "全部类型" ("F" column)
$allUpdates = @{
    13 = 7499
    14 = 7690
    15 = 24
    17 = 6273
    18 = 3288
    26 = 3683
    32 = 1326
    35 = 2643
    36 = 1574
    37 = 16
    40 = 3358
    41 = 200
    46 = 1313
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
